$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "3rd commit"
$ws.Range("A2").Select()
